$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 955.2381
$ws.Range("I43").Value = 1163.3334
$ws.Range("J43").Value = 872
$ws.Range("K43").Value = 1163.3334
$ws.Range("L43").Value = 872
$ws.Range("M43").Value = -1094.3334
$ws.Range("N43").Value = -1010
$ws.Range("H69").Value = 3384.2
$ws.Range("I69").Value = 2950.75
$ws.Range("J69").Value = 3541.818
$ws.Range("K69").Value = 8852.25
$ws.Range("L69").Value = 10625.454
$ws.Range("M69").Value = -7978.25
$ws.Range("N69").Value = -12373.454
$ws.Range("H72").Value = 3384.2
$ws.Range("I72").Value = 2950.75
$ws.Range("J72").Value = 3541.818
$ws.Range("K72").Value = 26556.75
$ws.Range("L72").Value = 31876.362
$ws.Range("M72").Value = -22188.75
$ws.Range("N72").Value = -40612.362
$ws.Range("H132").Value = 1423368.8
$ws.Range("I132").Value = 1595426
$ws.Range("K132").Value = 4786278
$ws.Range("M132").Value = -4783748
$ws.Range("H137").Value = 1217.2659
$ws.Range("I137").Value = 821.1667
$ws.Range("J137").Value = 1390.1091
$ws.Range("K137").Value = 2463.5001
$ws.Range("L137").Value = 4170.3273
$ws.Range("M137").Value = 86.4998999999998
$ws.Range("N137").Value = -9270.327300000001
$ws.Range("H138").Value = 6984.3076
$ws.Range("J138").Value = 8399.6
$ws.Range("L138").Value = 25198.8
$ws.Range("N138").Value = -35478.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 7208.3335
$ws.Range("I19").Value = 2500
$ws.Range("J19").Value = 8777.777
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 8777.777
$ws.Range("M19").Value = -2271
$ws.Range("N19").Value = -9235.777
$ws.Range("H61").Value = 1856.5676
$ws.Range("I61").Value = 1693.5312
$ws.Range("J61").Value = 2900
$ws.Range("K61").Value = 1693.5312
$ws.Range("L61").Value = 2900
$ws.Range("M61").Value = -1481.5312
$ws.Range("N61").Value = -3324
$ws.Range("H97").Value = 1505
$ws.Range("I97").Value = 1394.3334
$ws.Range("J97").Value = 2003
$ws.Range("K97").Value = 1394.3334
$ws.Range("L97").Value = 2003
$ws.Range("M97").Value = -898.3334
$ws.Range("N97").Value = -2995
$ws.Range("H102").Value = 1803.1025
$ws.Range("I102").Value = 1694.5714
$ws.Range("K102").Value = 1694.5714
$ws.Range("M102").Value = -72.57140000000004
$ws.Range("H132").Value = 1973.7819
$ws.Range("I132").Value = 1236.75
$ws.Range("J132").Value = 2999.2173
$ws.Range("K132").Value = 3710.25
$ws.Range("L132").Value = 8997.651899999999
$ws.Range("M132").Value = -1180.25
$ws.Range("N132").Value = -14057.6519
$ws.Range("H136").Value = 1856.5676
$ws.Range("I136").Value = 1693.5312
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 5080.5936
$ws.Range("L136").Value = 8700
$ws.Range("M136").Value = -2530.5936
$ws.Range("N136").Value = -13800
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2079
$ws.Range("I86").Value = 2095.4119
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 2095.4119
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -972.4119000000001
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 2079
$ws.Range("I89").Value = 2095.4119
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 10477.0595
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -4861.059499999999
$ws.Range("N89").Value = -20232
$ws.Range("H94").Value = 666.6667
$ws.Range("I94").Value = 760
$ws.Range("J94").Value = 200
$ws.Range("K94").Value = 760
$ws.Range("L94").Value = 200
$ws.Range("M94").Value = -309
$ws.Range("N94").Value = -1102
$ws.Range("H99").Value = 2192.5715
$ws.Range("I99").Value = 2399.6365
$ws.Range("J99").Value = 1433.3334
$ws.Range("K99").Value = 2399.6365
$ws.Range("L99").Value = 1433.3334
$ws.Range("M99").Value = -901.6365000000001
$ws.Range("N99").Value = -4429.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 36
$ws.Range("I7").Value = 38.375
$ws.Range("J7").Value = 34.1
$ws.Range("K7").Value = 38.375
$ws.Range("L7").Value = 34.1
$ws.Range("M7").Value = 74.625
$ws.Range("N7").Value = -260.1
$ws.Range("H16").Value = 3427.6667
$ws.Range("I16").Value = 2090.8462
$ws.Range("K16").Value = 2090.8462
$ws.Range("M16").Value = -1803.8462
$ws.Range("H105").Value = 881.3333
$ws.Range("I105").Value = 837
$ws.Range("J105").Value = 970
$ws.Range("K105").Value = 837
$ws.Range("L105").Value = 970
$ws.Range("M105").Value = 910
$ws.Range("N105").Value = -4464
$ws.Range("H113").Value = 3427.6667
$ws.Range("I113").Value = 2090.8462
$ws.Range("K113").Value = 2090.8462
$ws.Range("M113").Value = 79.15380000000005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 6271.4287
$ws.Range("I74").Value = 1100
$ws.Range("J74").Value = 8340
$ws.Range("K74").Value = 3300
$ws.Range("L74").Value = 25020
$ws.Range("M74").Value = -2239
$ws.Range("N74").Value = -27142
$ws.Range("H77").Value = 6271.4287
$ws.Range("I77").Value = 1100
$ws.Range("J77").Value = 8340
$ws.Range("K77").Value = 9900
$ws.Range("L77").Value = 75060
$ws.Range("M77").Value = -4596
$ws.Range("N77").Value = -85668
$ws.Range("H120").Value = 12000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 12000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 36000
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -45676
$ws.Range("H123").Value = 1065
$ws.Range("I123").Value = 230
$ws.Range("K123").Value = 690
$ws.Range("M123").Value = 1760
$ws.Range("H131").Value = 864.5268600000001
$ws.Range("I131").Value = 369.17648
$ws.Range("J131").Value = 975.3289
$ws.Range("K131").Value = 1107.52944
$ws.Range("L131").Value = 2925.9867
$ws.Range("M131").Value = 3932.47056
$ws.Range("N131").Value = -13005.9867
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5250.143
$ws.Range("I70").Value = 5120.35
$ws.Range("J70").Value = 5574.625
$ws.Range("K70").Value = 5120.35
$ws.Range("L70").Value = 5574.625
$ws.Range("M70").Value = -4850.35
$ws.Range("N70").Value = -6114.625
$ws.Range("H73").Value = 5250.143
$ws.Range("I73").Value = 5120.35
$ws.Range("J73").Value = 5574.625
$ws.Range("K73").Value = 5120.35
$ws.Range("L73").Value = 5574.625
$ws.Range("M73").Value = -4184.35
$ws.Range("N73").Value = -7446.625
$ws.Range("H96").Value = 18233.334
$ws.Range("J96").Value = 18233.334
$ws.Range("L96").Value = 18233.334
$ws.Range("N96").Value = -23725.334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2494.125
$ws.Range("H132").Value = 2856.1304
$ws.Range("I132").Value = 2100.3333
$ws.Range("J132").Value = 3680.6365
$ws.Range("K132").Value = 6300.999899999999
$ws.Range("L132").Value = 11041.9095
$ws.Range("M132").Value = -3770.999899999999
$ws.Range("N132").Value = -16101.9095
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 33274
$ws.Range("J95").Value = 33274
$ws.Range("L95").Value = 33274
$ws.Range("N95").Value = -38766
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -7340
